# Applies the "UUSI PYTHON" / "Uusi "-"Testaa" column edit described in the
# commit diff:
#   - Insert a new column before the existing "Turha"/"turhuus" column (C)
#     and fill it with the header "UUSI PYTHON" plus the value "UUSI PYTHON"
#     repeated for every data row.
#   - Append a brand-new trailing column with header "Uusi " and the value
#     "Testaa" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill what will become the new trailing column F ("Uusi " header /
#     "Testaa" data) first, while it's still column E, so the new
#     shared-string entries land in the same order as the original edit
#     (Uusi, Testaa, UUSI PYTHON). It shifts to F once column C is inserted
#     below.
$ws.Cells.Item(1, 5).Value = "Uusi "
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Value = "Testaa"
}
$ws.Range("E1:E6").Font.Size = 15

# --- Insert a new column C (pushes old "Turha"/"Tsekkaus"/"Uusi " columns
#     right by one, landing on D/E/F) ---
$ws.Columns.Item(3).Insert()

# Header + data for the newly inserted column C ("UUSI PYTHON")
$ws.Cells.Item(1, 3).Value = "UUSI PYTHON"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = "UUSI PYTHON"
}

# Match the look of the rest of the header/data (same font size as cols A/B)
$ws.Range("C1:C6").Font.Size = 15

# Give the new column a bit more width, like the real edit did.
$ws.Columns.Item(3).ColumnWidth = 18.26953125

# Update selection to mirror the saved state in the diff.
$ws.Range("C10").Select()
